$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 213 ("Fukumoto" row), shifting
# the existing rows 213:224 down to 215:226.
$ws.Rows("213:214").Insert()

# Row 213 - new Valencia "Primera" weekly record
$ws.Cells.Item(213, 1).Value = 11
$ws.Cells.Item(213, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(213, 3).Value = "Bíobío"
$ws.Cells.Item(213, 4).Value = 44610
$ws.Cells.Item(213, 5).Value = 8
$ws.Cells.Item(213, 6).Value = "Fruta"
$ws.Cells.Item(213, 7).Value = 100102
$ws.Cells.Item(213, 8).Value = "Cítricos"
$ws.Cells.Item(213, 9).Value = 100102005
$ws.Cells.Item(213, 10).Value = "Naranja"
$ws.Cells.Item(213, 11).Value = "Valencia"
$ws.Cells.Item(213, 12).Value = "Primera"
$ws.Cells.Item(213, 13).Value = 180
$ws.Cells.Item(213, 14).Value = 9500
$ws.Cells.Item(213, 15).Value = 10000
$ws.Cells.Item(213, 16).Value = 9722
$ws.Cells.Item(213, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(213, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(213, 19).Value = 648
$ws.Cells.Item(213, 20).Value = 15

# Row 214 - new Valencia "Segunda" weekly record
$ws.Cells.Item(214, 1).Value = 11
$ws.Cells.Item(214, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(214, 3).Value = "Bíobío"
$ws.Cells.Item(214, 4).Value = 44610
$ws.Cells.Item(214, 5).Value = 8
$ws.Cells.Item(214, 6).Value = "Fruta"
$ws.Cells.Item(214, 7).Value = 100102
$ws.Cells.Item(214, 8).Value = "Cítricos"
$ws.Cells.Item(214, 9).Value = 100102005
$ws.Cells.Item(214, 10).Value = "Naranja"
$ws.Cells.Item(214, 11).Value = "Valencia"
$ws.Cells.Item(214, 12).Value = "Segunda"
$ws.Cells.Item(214, 13).Value = 160
$ws.Cells.Item(214, 14).Value = 8000
$ws.Cells.Item(214, 15).Value = 8500
$ws.Cells.Item(214, 16).Value = 8312
$ws.Cells.Item(214, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(214, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(214, 19).Value = 554
$ws.Cells.Item(214, 20).Value = 15

# Ensure the date cells keep the date-number-format style used by column D.
$ws.Range("D213:D214").NumberFormat = $ws.Range("D215").NumberFormat
